# Applies the cryptos.xlsx price/volume update diff via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '90.430.48'
$ws.Cells.Item(2, 5).Value = '  -0.53%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.059.17'

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '242.28'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.98%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '615.09'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -2.46%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +7.19%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.362'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.52%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.03%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.056.41'
$ws.Cells.Item(10, 5).Value = '  -1.71%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.730'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +2.43%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.44%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000244'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.17%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '34.54'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -5.25%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '90.466.87'
$ws.Cells.Item(15, 5).Value = '  +0.06%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.41'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.20%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.622.71'
$ws.Cells.Item(17, 5).Value = '  -1.75%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.020.51'
$ws.Cells.Item(18, 5).Value = '  -2.45%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.61'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -3.02%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '14.28'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.48%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000208'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.16%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.71'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +3.61%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '437.68'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.38%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.93'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.36%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'NEARProtocol'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.55'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.63%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Litecoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '90.20'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.97%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '11.67'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -6.29%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '3.246.60'
$ws.Cells.Item(28, 5).Value = '  -0.85%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.03%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Cronos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.180'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +13.02%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.240'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +25.83%  '

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +12.83%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.01'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -4.93%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +31.28%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.166'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +11.53%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.57'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +7.02%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.13'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.21%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.15'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +30.16%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.89'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.88%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '484.69'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -4.64%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.48'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -7.73%  '

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.27'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -1.05%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.14'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.27%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.01%  '

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '153.89'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.56%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.86'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -2.08%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.673'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.62%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '44.10'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.71%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.31'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.58%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.37'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.59%  '
